# Re-run of the depth-multiplier model: rows 105-143 (columns C:J) are
# refreshed with the latest model outputs. Columns A and B are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colVals = @{
    3  = 0.0187499999985
    4  = 0.12946428575
    5  = 0.2633928571249999
    6  = 0.3973214285
    7  = 0.53125
    8  = 0.6651785714999999
    9  = 0.7991071428749998
    10 = 0.93303571425
}

for ($r = 105; $r -le 143; $r++) {
    foreach ($col in $colVals.Keys) {
        $ws.Cells.Item($r, $col).Value = $colVals[$col]
    }
}
